# Add a new worksheet ("Sheet27") that duplicates the last sheet ("Sheet26")
# of the workbook, but with the "author" column (column B) rewritten from
# "Maja Založnik" to the newly-inserted author "Marko" (pulled from the db).

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Sheet26")

# Mark the source sheet as fully selected (closest reachable approximation of
# the author having select-alled the sheet before branching off a copy).
$src.Activate()
$src.Range("H10").Select() | Out-Null
$src.Cells.Select() | Out-Null

# Duplicate Sheet26 right after itself, then rename the copy to Sheet27.
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "Sheet27"

# Swap the author for the freshly inserted rows: "Marko" replaces
# "Maja Založnik" in column B (author) for every data row.
$newSheet.Range("B2:B5").Value = "Marko"

# Leave the new sheet active/selected, cursor resting just below the data.
$newSheet.Activate()
$newSheet.Range("B6").Select() | Out-Null
